$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cryptocurrency Price (D) and Volume(1h) (E) cells with refreshed
# values from the latest data pull. Values are plain text in the sheet, so
# for any new value that looks like a numeric literal we force text storage
# (leading apostrophe) and then reset the cell style back to Normal so no
# numeric auto-conversion or stray formatting is introduced.

$ws.Range('D2').Value = '26.671.51'
$ws.Range('D3').Value = '1.599.33'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.34%  '
$cell = $ws.Range('D5')
$cell.Value = "'211.66"
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  +0.32%  '
$cell = $ws.Range('D8')
$cell.Value = "'0.0619"
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -0.10%  '
$cell = $ws.Range('D9')
$cell.Value = "'0.247"
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.12%  '
$cell = $ws.Range('D10')
$cell.Value = "'19.61"
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -0.91%  '
$cell = $ws.Range('D11')
$cell.Value = "'0.0840"
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '1.823.17'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '1.599.69'
$ws.Range('E13').Value = '  -0.07%  '
$cell = $ws.Range('D14')
$cell.Value = "'4.03"
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('E15').Value = '  -0.02%  '
$cell = $ws.Range('D16')
$cell.Value = "'65.02"
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '26.664.15'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '0.0₃0736'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('E19').Value = '  +0.29%  '
$cell = $ws.Range('D20')
$cell.Value = "'208.35"
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -1.01%  '
$cell = $ws.Range('D21')
$cell.Value = "'7.06"
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +5.06%  '
$cell = $ws.Range('D22')
$cell.Value = "'4.29"
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('E23').Value = '  +0.48%  '
$cell = $ws.Range('D24')
$cell.Value = "'8.96"
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +0.46%  '
$cell = $ws.Range('D25')
$cell.Value = "'145.14"
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -1.06%  '
$ws.Range('E26').Value = '  +0.41%  '
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('E29').Value = '  -0.37%  '
$cell = $ws.Range('D30')
$cell.Value = "'0.0513"
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +2.00%  '
$cell = $ws.Range('D31')
$cell.Value = "'1.16"
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('E33').Value = '  +1.45%  '
$ws.Range('D34').Value = '1.277.93'
$ws.Range('E34').Value = '  -1.65%  '
$cell = $ws.Range('D35')
$cell.Value = "'0.618"
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -8.04%  '
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('E38').Value = '  -0.99%  '
$cell = $ws.Range('D39')
$cell.Value = "'0.836"
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('E40').Value = '  +15.84%  '
$ws.Range('E41').Value = '  +2.80%  '
$ws.Range('E42').Value = '  +0.10%  '
$cell = $ws.Range('D43')
$cell.Value = "'0.784"
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -0.73%  '
$cell = $ws.Range('D44')
$cell.Value = "'64.03"
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').Value = '1.735.49'
$ws.Range('E45').Value = '  +0.03%  '
$cell = $ws.Range('D46')
$cell.Value = "'91.00"
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('E47').Value = '  -2.43%  '
$ws.Range('D48').Value = '0.0₆0105'
$ws.Range('E48').Value = '  +0.05%  '
$cell = $ws.Range('D49')
$cell.Value = "'0.102"
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('E51').Value = '  +0.05%  '
